# Scheduled-runner market price refresh: recompute currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ (columns H-N) for the affected Leve rows across
# the per-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW).
#
# Some rows lose their HQ-profit (column N) or NQ-profit (column M) cell
# entirely once the corresponding price data is no longer meaningful (e.g.
# HQ price drops to 0), while others gain a profit cell that didn't exist
# before. ClearContents() is used so such cells are removed outright rather
# than merely blanked.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H74").Value = 3288.75
$ws.Range("I74").Value = 2983.9
$ws.Range("J74").Value = 3593.6
$ws.Range("K74").Value = 2983.9
$ws.Range("L74").Value = 3593.6
$ws.Range("M74").Value = -2047.9
$ws.Range("N74").Value = -5465.6

$ws.Range("H77").Value = 3288.75
$ws.Range("I77").Value = 2983.9
$ws.Range("J77").Value = 3593.6
$ws.Range("K77").Value = 14919.5
$ws.Range("L77").Value = 17968
$ws.Range("M77").Value = -10239.5
$ws.Range("N77").Value = -27328

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 852
$ws.Range("I2").Value = 841.4286
$ws.Range("K2").Value = 841.4286
$ws.Range("M2").Value = -728.4286

$ws.Range("H32").Value = 14088177
$ws.Range("I32").Value = 3560.2131
$ws.Range("K32").Value = 3560.2131
$ws.Range("M32").Value = -3273.2131

$ws.Range("H61").Value = 6411559
$ws.Range("I61").Value = 6411559
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6411559
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6411347
$ws.Range("N61").ClearContents()

$ws.Range("H97").Value = 1365
$ws.Range("I97").Value = 1248
$ws.Range("J97").Value = 1950
$ws.Range("K97").Value = 1248
$ws.Range("L97").Value = 1950
$ws.Range("M97").Value = -752
$ws.Range("N97").Value = -2942

$ws.Range("H116").Value = 852
$ws.Range("I116").Value = 841.4286
$ws.Range("K116").Value = 841.4286
$ws.Range("M116").Value = 1452.5714

$ws.Range("H132").Value = 1085.8695
$ws.Range("I132").Value = 768.2564
$ws.Range("J132").Value = 2855.4285
$ws.Range("K132").Value = 2304.7692
$ws.Range("L132").Value = 8566.2855
$ws.Range("M132").Value = 225.2308000000003
$ws.Range("N132").Value = -13626.2855

$ws.Range("H136").Value = 6411559
$ws.Range("I136").Value = 6411559
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 19234677
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -19232127
$ws.Range("N136").ClearContents()

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 852
$ws.Range("I3").Value = 841.4286
$ws.Range("K3").Value = 841.4286
$ws.Range("M3").Value = -727.4286

$ws.Range("H134").Value = 2646491.5
$ws.Range("I134").Value = 898.2564
$ws.Range("J134").Value = 37039204
$ws.Range("K134").Value = 2694.7692
$ws.Range("L134").Value = 111117612
$ws.Range("M134").Value = -159.7691999999997
$ws.Range("N134").Value = -111122682

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 1544119.9
$ws.Range("I31").Value = 1792936
$ws.Range("J31").Value = 1460
$ws.Range("K31").Value = 1792936
$ws.Range("L31").Value = 1460
$ws.Range("M31").Value = -1792641
$ws.Range("N31").Value = -2050

$ws.Range("H34").Value = 1544119.9
$ws.Range("I34").Value = 1792936
$ws.Range("J34").Value = 1460
$ws.Range("K34").Value = 1792936
$ws.Range("L34").Value = 1460
$ws.Range("M34").Value = -1792734
$ws.Range("N34").Value = -1864

$ws.Range("H94").Value = 13514.25
$ws.Range("I94").Value = 34000
$ws.Range("J94").Value = 1222.8
$ws.Range("K94").Value = 34000
$ws.Range("L94").Value = 1222.8
$ws.Range("M94").Value = -33549
$ws.Range("N94").Value = -2124.8

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H63").Value = 1058.3334
$ws.Range("I63").Value = 687.5
$ws.Range("J63").Value = 1800
$ws.Range("K63").Value = 2062.5
$ws.Range("L63").Value = 5400
$ws.Range("M63").Value = -1313.5
$ws.Range("N63").Value = -6898

$ws.Range("H64").Value = 2000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 6000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -6540

$ws.Range("H66").Value = 1058.3334
$ws.Range("I66").Value = 687.5
$ws.Range("J66").Value = 1800
$ws.Range("K66").Value = 6187.5
$ws.Range("L66").Value = 16200
$ws.Range("M66").Value = -2443.5
$ws.Range("N66").Value = -23688

$ws.Range("H67").Value = 2000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 6000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -7872

$ws.Range("H87").Value = 7671.3335
$ws.Range("I87").Value = 7014
$ws.Range("K87").Value = 21042
$ws.Range("M87").Value = -19794

$ws.Range("H90").Value = 7671.3335
$ws.Range("I90").Value = 7014
$ws.Range("K90").Value = 63126
$ws.Range("M90").Value = -56886

$ws.Range("H114").Value = 655.4838999999999
$ws.Range("I114").Value = 259.5
$ws.Range("J114").Value = 844.0476
$ws.Range("K114").Value = 778.5
$ws.Range("L114").Value = 2532.1428
$ws.Range("M114").Value = 2475.5
$ws.Range("N114").Value = -9040.1428

$ws.Range("H117").Value = 1217.4546
$ws.Range("I117").Value = 352
$ws.Range("J117").Value = 1938.6666
$ws.Range("K117").Value = 1056
$ws.Range("L117").Value = 5815.9998
$ws.Range("M117").Value = 2386
$ws.Range("N117").Value = -12699.9998

$ws.Range("H129").Value = 1429.4286
$ws.Range("I129").Value = 1135
$ws.Range("J129").Value = 1547.2
$ws.Range("K129").Value = 3405
$ws.Range("L129").Value = 4641.6
$ws.Range("M129").Value = 1595
$ws.Range("N129").Value = -14641.6

$ws.Range("H131").Value = 784.1
$ws.Range("J131").Value = 818.6517
$ws.Range("L131").Value = 2455.9551
$ws.Range("N131").Value = -12535.9551

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H122").Value = 100021180
$ws.Range("I122").Value = 100021180
$ws.Range("K122").Value = 300063540
$ws.Range("M122").Value = -300061090

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H61").Value = 2045.1111
$ws.Range("I61").Value = 2045.1111
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2045.1111
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1843.1111
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 2045.1111
$ws.Range("I113").Value = 2045.1111
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2045.1111
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 124.8888999999999
$ws.Range("N113").ClearContents()
